$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.064.00"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "'2.021.06"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'225.95"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "'0.608"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'54.81"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").Value = "'0.0782"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("E11").Value = "  -4.66%  "
$ws.Range("D12").Value = "'2.322.67"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "'14.08"
$ws.Range("E13").Value = "  -4.70%  "
$ws.Range("D14").Value = "'20.07"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").Value = "'0.739"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").Value = "'5.19"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").Value = "'2.095.56"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "'36.978.59"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").Value = "'6.43"
$ws.Range("E19").Value = "  +5.48%  "
$ws.Range("D20").Value = "'68.70"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "'0.0₃0813"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").Value = "'222.59"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").Value = "'2.17"
$ws.Range("E25").Value = "  -5.77%  "
$ws.Range("D26").Value = "'165.07"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").Value = "'9.14"
$ws.Range("E27").Value = "  -5.52%  "
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("D29").Value = "'18.59"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "'4.47"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").Value = "'0.0602"
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("D35").Value = "'2.33"
$ws.Range("E35").Value = "  -3.97%  "
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'3.08"
$ws.Range("E38").Value = "  -5.06%  "
$ws.Range("D39").Value = "'5.37"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").Value = "'1.458.52"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").Value = "'0.0212"
$ws.Range("E41").Value = "  -4.61%  "
$ws.Range("D42").Value = "'95.22"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("D43").Value = "'2.80"
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").Value = "'0.0909"
$ws.Range("E44").Value = "  -3.34%  "
$ws.Range("D45").Value = "'16.17"
$ws.Range("E45").Value = "  -6.15%  "
$ws.Range("E46").Value = "  -3.17%  "
$ws.Range("D47").Value = "'7.18"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").Value = "'2.92"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'2.211.76"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -9.19%  "
